{"js": "// Fix a double space typo in the project-description run:\n// \"controle de  luminosidade\" -> \"controle de luminosidade\"\nconst searchText = \"controle de  luminosidade\";\nconst replacement = \"controle de luminosidade\";\n\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(replacement, \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix a double space typo in the project-description run:\n# \"controle de  luminosidade\" -> \"controle de luminosidade\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"controle de  luminosidade\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"controle de luminosidade\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$find.MatchWildcards, $null, $null, [ref]$find.Forward, $null, $null, [ref]$find.Replacement.Text, 2) | Out-Null\n"}
